# "Generate Report for Handback"
# Updates the generated timestamp values on the Overview/zh-cn/de-de sheets
# to reflect a fresh handback report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G)
$wsOverview.Range("G2").Value = "2016-08-27 09:04:17"

# zh-cn sheet: "Correspond Handoff Datetime" (column H) and
# "Correspond Handback DateTime" (column K)
$wsZhCn.Range("H2").Value = "2016-08-27 09:04:12"
$wsZhCn.Range("K2").Value = "2016-08-27 09:04:29"

# de-de sheet: "Correspond Handoff Datetime" (column H) and
# "Correspond Handback DateTime" (column K)
$wsDeDe.Range("H2").Value = "2016-08-27 09:04:17"
$wsDeDe.Range("K2").Value = "2016-08-27 09:04:35"
